$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the existing header style (bold, centered, bordered) used by B1:H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-22 for columns I and J
$values = @(
    @(9,9),
    @(7,7),
    @(8,8),
    @(9,9),
    @(7,7),
    @(9,9),
    @(7,7),
    @(9,9),
    @(8,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(9,9),
    @(8,8),
    @(7,7),
    @(4,5),
    @(6,6),
    @(5,5),
    @(5,5),
    @(4,4),
    @(9,9)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
